# Weekly refresh of the Papaya price table (Vega Central Mapocho de Santiago).
# Each data row (2-31) keeps its place in the sheet, but its Fecha / Calidad /
# Volumen / Precio / Unidad fields are updated to the new weekly figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44880
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("S2").Value = 2000
$ws.Range("D3").Value = 44880
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 180
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 1500
$ws.Range("D4").Value = 44292
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 1400
$ws.Range("D5").Value = 44343
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 1700
$ws.Range("O5").Value = 1700
$ws.Range("P5").Value = 1700
$ws.Range("Q5").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S5").Value = 1700
$ws.Range("T5").Value = 1
$ws.Range("D6").Value = 45113
$ws.Range("M6").Value = 6
$ws.Range("D7").Value = 45113
$ws.Range("M7").Value = 8
$ws.Range("D8").Value = 45113
$ws.Range("M8").Value = 15
$ws.Range("D9").Value = 45113
$ws.Range("M9").Value = 8
$ws.Range("D10").Value = 45126
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 55
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("S10").Value = 2400
$ws.Range("D11").Value = 45126
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = '$/bandeja 10 kilos'
$ws.Range("S11").Value = 2000
$ws.Range("T11").Value = 10
$ws.Range("D12").Value = 45126
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/bandeja 10 kilos'
$ws.Range("S12").Value = 1500
$ws.Range("T12").Value = 10
$ws.Range("D13").Value = 45126
$ws.Range("L13").Value = 'Tercera'
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 1200
$ws.Range("D14").Value = 45125
$ws.Range("L14").Value = 'Especial'
$ws.Range("M14").Value = 55
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 24000
$ws.Range("P14").Value = 24000
$ws.Range("S14").Value = 2400
$ws.Range("D15").Value = 45125
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("Q15").Value = '$/bandeja 10 kilos'
$ws.Range("S15").Value = 2000
$ws.Range("T15").Value = 10
$ws.Range("D16").Value = 45125
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 1500
$ws.Range("D17").Value = 45125
$ws.Range("L17").Value = 'Tercera'
$ws.Range("M17").Value = 45
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("S17").Value = 1200
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 140
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("S18").Value = 2400
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("S19").Value = 2000
$ws.Range("D20").Value = 45118
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 180
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 15000
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("T20").Value = 10
$ws.Range("D21").Value = 45118
$ws.Range("L21").Value = 'Tercera'
$ws.Range("M21").Value = 75
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/bandeja 10 kilos'
$ws.Range("S21").Value = 1200
$ws.Range("T21").Value = 10
$ws.Range("D22").Value = 44371
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 1800
$ws.Range("O22").Value = 1800
$ws.Range("P22").Value = 1800
$ws.Range("S22").Value = 1800
$ws.Range("D23").Value = 44371
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 1200
$ws.Range("O23").Value = 1200
$ws.Range("P23").Value = 1200
$ws.Range("S23").Value = 1200
$ws.Range("D24").Value = 44904
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 45
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("Q24").Value = '$/bandeja 10 kilos'
$ws.Range("S24").Value = 1500
$ws.Range("T24").Value = 10
$ws.Range("D25").Value = 44904
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("S25").Value = 1000
$ws.Range("D26").Value = 44391
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 15
$ws.Range("N26").Value = 1500
$ws.Range("O26").Value = 1500
$ws.Range("P26").Value = 1500
$ws.Range("Q26").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S26").Value = 1500
$ws.Range("T26").Value = 1
$ws.Range("D27").Value = 44391
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 1000
$ws.Range("O27").Value = 1000
$ws.Range("P27").Value = 1000
$ws.Range("Q27").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S27").Value = 1000
$ws.Range("T27").Value = 1
$ws.Range("D28").Value = 44336
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 10
$ws.Range("N28").Value = 1500
$ws.Range("O28").Value = 1500
$ws.Range("P28").Value = 1500
$ws.Range("Q28").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("T28").Value = 1
$ws.Range("D29").Value = 44195
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 15000
$ws.Range("S29").Value = 1500
$ws.Range("D30").Value = 44400
$ws.Range("M30").Value = 25
$ws.Range("N30").Value = 1500
$ws.Range("O30").Value = 1500
$ws.Range("P30").Value = 1500
$ws.Range("Q30").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S30").Value = 1500
$ws.Range("T30").Value = 1
$ws.Range("D31").Value = 44309
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 10
$ws.Range("N31").Value = 1600
$ws.Range("O31").Value = 1600
$ws.Range("P31").Value = 1600
$ws.Range("Q31").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S31").Value = 1600
$ws.Range("T31").Value = 1
